$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generations 0-17 (rows 2-19) -> Fitness value 7900
$ws.Range("C2:C19").Value = 7900

# Generations 18-250 (rows 20-252) -> Fitness value 7293
$ws.Range("C20:C252").Value = 7293
